# Insert a new data row at row 548 (pushes existing rows 548..616 down to 549..617)
# and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(548).Insert()

$ws.Cells.Item(548, 1).Value  = 3
$ws.Cells.Item(548, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(548, 3).Value  = "Coquimbo"
$ws.Cells.Item(548, 4).Value  = 45142
$ws.Cells.Item(548, 5).Value  = 5
$ws.Cells.Item(548, 6).Value  = 100112040
$ws.Cells.Item(548, 7).Value  = "Cilantro"
$ws.Cells.Item(548, 8).Value  = "Sin especificar"
$ws.Cells.Item(548, 9).Value  = "Primera"
$ws.Cells.Item(548, 10).Value = 120
$ws.Cells.Item(548, 11).Value = 4000
$ws.Cells.Item(548, 12).Value = 4000
$ws.Cells.Item(548, 13).Value = 4000
$ws.Cells.Item(548, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(548, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(548, 16).Value = 1333
$ws.Cells.Item(548, 17).Value = 3
$ws.Cells.Item(548, 18).Value = "Hortaliza"
